# Refresh the cryptos price/volume table (and fix the TRON / WrappedBTC row
# order, which the source swapped back) to match the latest scrape.
#
# NOTE: several Price values look like plain numbers (e.g. "1.00", "0.121").
# Assigning such a string straight to .Value lets Excel's smart typing coerce
# it into a numeric cell (dropping the significant trailing zeros / decimal
# precision the source text needs). To keep those cells as literal text we
# temporarily force a Text number format before the write, then restore the
# default "Normal" style afterwards so the cell's formatting/style index is
# unaffected by the round trip.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.296.79"
$ws.Range("E2").Value = "  +7.20%  "
$ws.Range("D3").Value = "3.570.86"
$ws.Range("E3").Value = "  +11.15%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "188.06"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +9.78%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "553.32"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +5.77%  "
$ws.Range("D7").Value = "3.563.56"
$ws.Range("E7").Value = "  +11.12%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.608"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.81%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.635"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +5.48%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.153"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +15.90%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "54.85"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +4.15%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000271"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +8.22%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.43"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.19%  "
$ws.Range("D15").Value = "4.132.77"
$ws.Range("E15").Value = "  +10.95%  "
$ws.Range("D16").Value = "3.571.07"
$ws.Range("E16").Value = "  +11.07%  "
$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "67.438.62"
$ws.Range("E17").Value = "  +7.57%  "
$ws.Range("B18").Value = "TRON"
$ws.Range("C18").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.121"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +5.21%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "18.26"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +6.80%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.00"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +9.37%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.995"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.48%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "431.65"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +18.74%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "85.67"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +5.82%  "
$ws.Range("E24").Value = "  +5.04%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.13"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +5.67%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.12"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.27%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.90"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +10.86%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.17"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.84%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "12.09"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +7.99%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.02"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +11.39%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "30.55"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +7.90%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "644.59"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.59%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.66"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.35%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "11.75"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +5.11%  "
$ws.Range("E35").Value = "  +6.36%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "59.77"
$ws.Range("D36").Style = "Normal"
$ws.Range("E37").Value = "  +25.23%  "
$ws.Range("D38").Value = "0.0₃0825"
$ws.Range("E38").Value = "  +17.53%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "38.56"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +5.74%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.999"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.16%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.391"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.08%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.35"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +15.65%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.00"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.03%  "
$ws.Range("D44").Value = "3.054.30"
$ws.Range("E44").Value = "  +6.75%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.66"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +5.17%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.43"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +12.76%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.88"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +12.43%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.83"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +5.21%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0419"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +7.68%  "
$ws.Range("E50").Value = "  +5.44%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.75"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +14.07%  "
